$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price (column D) and 1h volume change (column E)
# values for the rows that changed in this data refresh.
# A leading apostrophe forces Excel to keep plain-number-looking price
# strings (e.g. "305.42") as text, matching the original inlineStr cells.

$ws.Range('D2').Value = '42.643.56'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '2.280.77'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''305.42'
$ws.Range('D6').Value = '''96.43'
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -3.43%  '
$ws.Range('D10').Value = '''35.50'
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '''18.29'
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('E14').Value = '  -2.27%  '
$ws.Range('D15').Value = '2.635.75'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '2.284.36'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('D17').Value = '''0.778'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '42.580.28'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '''12.92'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').Value = '''67.11'
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').Value = '''235.84'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  -2.71%  '
$ws.Range('D25').Value = '''2.45'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '''4.02'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '''25.11'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').Value = '''166.18'
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('D32').Value = '''33.04'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').Value = '''4.75'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('D35').Value = '''4.96'
$ws.Range('E35').Value = '  -3.29%  '
$ws.Range('D36').Value = '''17.58'
$ws.Range('E36').Value = '  -3.41%  '
$ws.Range('D37').Value = '''2.39'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('D41').Value = '''0.109'
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').Value = '2.000.15'
$ws.Range('E43').Value = '  -0.72%  '
$ws.Range('E44').Value = '  -2.74%  '
$ws.Range('D45').Value = '''18.16'
$ws.Range('E45').Value = '  +3.65%  '
$ws.Range('D46').Value = '''9.97'
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('E47').Value = '  -8.35%  '
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('D49').Value = '''2.86'
$ws.Range('E49').Value = '  +4.95%  '
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('D51').Value = '2.504.15'
$ws.Range('E51').Value = '  -0.82%  '
